# Adds two new columns, "I0" (col I) and "IF" (col J), to the sheet,
# matching header styling from the existing "IP" header (col H),
# and fills in the per-row numeric values for rows 2-68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
$ws.Cells.Item(1, 9).Value  = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Copy the formatting of the existing "IP" header cell (H1) onto the
# two new header cells (I1:J1) so they pick up the same bold/border/
# alignment style used by the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-68): column I = "I0", column J = "IF" ---
$values = @(
    @(9,9),
    @(6,6),
    @(8,9),
    @(3,4),
    @(7,7),
    @(7,8),
    @(6,7),
    @(8,8),
    @(8,8),
    @(6,6),
    @(11,11),
    @(6,6),
    @(7,7),
    @(3,4),
    @(9,9),
    @(6,7),
    @(8,8),
    @(5,6),
    @(5,5),
    @(9,9),
    @(6,6),
    @(8,8),
    @(7,7),
    @(9,9),
    @(8,9),
    @(10,10),
    @(8,8),
    @(6,7),
    @(8,8),
    @(8,8),
    @(7,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,10),
    @(9,9),
    @(6,6),
    @(7,7),
    @(7,7),
    @(10,10),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(10,10),
    @(7,7),
    @(7,7),
    @(7,7),
    @(10,10),
    @(7,8),
    @(9,9),
    @(8,8),
    @(7,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(8,9),
    @(5,5),
    @(8,8),
    @(7,7),
    @(2,2)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value  = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}

Write-Host "Added I0/IF columns with $($values.Count) data rows"
